# Applies the edit described by the commit:
#   "Vorlage Aspektbearbeitung, Statusbericht, diverse Updates"
#
# Net effect inside word/document.xml:
#   1. After the "... der Arbeit:" run in the first paragraph, two new
#      runs are appended: a big-font space (title-style leftover
#      formatting) and the literal project code "MT-FS13-05".
#   2. A fresh, zero-width "_GoBack" bookmark is left right after that
#      new text (Word always drops one at the last edit position).
#   3. The previous "_GoBack" bookmark (that used to sit after the
#      "Master Thesis," run, further down in the document) is gone -
#      Word only ever keeps the most recent one.
#   4. Because a bookmark got removed and another got added, the
#      bookmark id Word assigns to "DDE_LINK" shifts from 0 to 1.

$d = $word.ActiveDocument

# --- remove the old "_GoBack" bookmark (Word keeps only the latest one) ---
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# --- locate the end of "... der Arbeit:" in the first paragraph ---
$p1 = $d.Paragraphs(1).Range
$r = $p1.Duplicate
$r.Find.Execute("der Arbeit:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

# Run 1: a single space, carrying the big leftover title-style formatting
# (as if the user had clicked at the end of the heading and started
# typing, inheriting the paragraph mark's run formatting, then fixed it
# up afterwards).
$r.InsertAfter(" ")
$r.Font.Name = "Constantia"
$r.Font.Bold = $false
$r.Font.BoldBi = $false
$r.Font.TextColor.ObjectThemeColor = 13
$r.Font.Kerning = 12
$r.Font.Size = 43

# Run 2: the actual project code, back in the paragraph's normal size.
$r.Collapse(0)
$r.InsertAfter("MT-FS13-05")
$r.Font.Bold = $false
$r.Font.Size = 12

# --- leave a new, zero-width "_GoBack" bookmark right after it ---
# (Insert a throwaway character first so the collapsed range sits away
# from the paragraph-mark boundary, add the bookmark there, then remove
# the throwaway character again - this keeps the bookmark anchored
# exactly between "MT-FS13-05" and the paragraph end.)
$r.Collapse(0)
$r.InsertAfter("Z")
$bmRange = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
$r.Text = ""
